# Create front-end menu for user interface
# Adds a new "status" lookup sheet (status_id / status / description) and
# wires it up to the "incident" sheet via a new status_id foreign-key column.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "status" worksheet as the last tab --------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$status = $wb.Worksheets.Add($null, $lastSheet)
$status.Name = "status"

# Fill column by column (A, then B, then C) so new shared strings are
# interned in the same order the original authoring session produced them.
$status.Range("A1").Value = "status_id"
$status.Range("A2").Value = 1
$status.Range("A3").Value = 2
$status.Range("A4").Value = 3

$status.Range("B1").Value = "status"
$status.Range("B2").Value = "Complete"
$status.Range("B3").Value = "In-Progress"
$status.Range("B4").Value = "Unsure"

$status.Range("C1").Value = "description"
$status.Range("C2").Value = "Incident resolved"
$status.Range("C3").Value = "Incident ongoing"
$status.Range("C4").Value = "Unclear as to resolution/continuation of incident"

# --- 2. Add the status_id foreign-key column to "incident" -----------------
$incident = $wb.Worksheets.Item("incident")

$incident.Range("D1").Value = "status"
$incident.Range("D2").Value = 1
$incident.Range("D3").Value = 2
$incident.Range("D4").Value = 3
$incident.Range("D5").Value = 1
$incident.Range("D6").Value = 2

# Select the populated range on the incident sheet (matches the selection
# state left behind after reviewing the new column).
$incident.Range("A1:D6").Select()

# --- 3. Leave the new "status" sheet active/selected, as the final tab ----
$status.Select()
$status.Range("E9").Select()
